$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Status: draft -> active
$ws.Range("B6").Value = "active"

# Date: bump to new publication timestamp
$ws.Range("B8").Value = "2024-12-16T14:50:05-03:00"

# Case Sensitive: false -> true
# (Assigning the literal string "true"/"false" directly would be auto-coerced
#  to a native boolean by Excel's input parser, which does not match the
#  original file's encoding of this value as plain text. Route the literal
#  through a formula and convert it to a static value via copy/paste so the
#  cell keeps its original "General" text styling and becomes a real text
#  value of "true" instead of a TRUE boolean.)
$c = $ws.Range("B17")
$c.Formula = "=""true"""
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0
